$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "Save" column, matching style of existing header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Fill H2:H34 with an indicator derived from column G (sum): 1 if sum > 10 else 0
$lastRow = 34
for ($r = 2; $r -le $lastRow; $r++) {
    $gVal = $ws.Cells.Item($r, 7).Value2
    if ($gVal -gt 10) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
